# Actualización automática de grupos experimentales
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap / set the Grupo_Experimental (column B) values for specific rows
$ws.Range("B2").Value = "Sin SmartScore"
$ws.Range("B4").Value = "Con SmartScore"
$ws.Range("B12").Value = "Sin SmartScore"
$ws.Range("B14").Value = "Con SmartScore"
$ws.Range("B15").Value = "Con SmartScore"
$ws.Range("B16").Value = "Sin SmartScore"
$ws.Range("B19").Value = "Con SmartScore"
$ws.Range("B20").Value = "Sin SmartScore"

# Convert SmartScore cells in row 20 from text to real numeric values
$ws.Range("I20").Value = 0.646
$ws.Range("L20").Value = 0.543
$ws.Range("O20").Value = 0.535
$ws.Range("R20").Value = 0.662
$ws.Range("U20").Value = 0.507
$ws.Range("X20").Value = 0.456
$ws.Range("AA20").Value = 0.72
$ws.Range("AD20").Value = 0.498
$ws.Range("AG20").Value = 0.472
